# Auto-generated: apply cached-value corrections for the scheduled market-data refresh
# (columns H-N hold scraped price snapshots / computed profit figures per sheet row; no formulas involved)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 8045.3
$ws.Range("I62").Value = 8106.5557
$ws.Range("K62").Value = 8106.5557
$ws.Range("M62").Value = -7482.5557
# Row 64
$ws.Range("H64").Value = 6356.8945
$ws.Range("I64").Value = 5000.091
$ws.Range("K64").Value = 5000.091
$ws.Range("M64").Value = -4752.091
# Row 65
$ws.Range("H65").Value = 8045.3
$ws.Range("I65").Value = 8106.5557
$ws.Range("K65").Value = 40532.7785
$ws.Range("M65").Value = -37412.7785
# Row 67
$ws.Range("H67").Value = 6356.8945
$ws.Range("I67").Value = 5000.091
$ws.Range("K67").Value = 5000.091
$ws.Range("M67").Value = -4142.091
# Row 80
$ws.Range("H80").Value = 706.5454999999999
$ws.Range("I80").Value = 654.9286
$ws.Range("J80").Value = 796.875
$ws.Range("K80").Value = 1964.7858
$ws.Range("L80").Value = 2390.625
$ws.Range("M80").Value = -966.7857999999999
$ws.Range("N80").Value = -4386.625
# Row 83
$ws.Range("H83").Value = 706.5454999999999
$ws.Range("I83").Value = 654.9286
$ws.Range("J83").Value = 796.875
$ws.Range("K83").Value = 5894.3574
$ws.Range("L83").Value = 7171.875
$ws.Range("M83").Value = -902.3573999999999
$ws.Range("N83").Value = -17155.875
# Row 137
$ws.Range("H137").Value = 35720676
$ws.Range("I137").Value = 55563588
$ws.Range("K137").Value = 166690764
$ws.Range("M137").Value = -166688214
# Row 138
$ws.Range("H138").Value = 8403.666999999999
$ws.Range("I138").Value = 5278.3335
$ws.Range("J138").Value = 8612.022000000001
$ws.Range("K138").Value = 15835.0005
$ws.Range("L138").Value = 25836.066
$ws.Range("M138").Value = -10695.0005
$ws.Range("N138").Value = -36116.06600000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1297.9697
$ws.Range("I2").Value = 1318.3182
$ws.Range("K2").Value = 1318.3182
$ws.Range("M2").Value = -1205.3182
# Row 32
$ws.Range("H32").Value = 19668.8
$ws.Range("I32").Value = 8213.444
$ws.Range("J32").Value = 43460.69
$ws.Range("K32").Value = 8213.444
$ws.Range("L32").Value = 43460.69
$ws.Range("M32").Value = -7926.444
$ws.Range("N32").Value = -44034.69
# Row 61
$ws.Range("H61").Value = 4686.8213
$ws.Range("I61").Value = 4049.76
$ws.Range("K61").Value = 4049.76
$ws.Range("M61").Value = -3837.76
# Row 63
$ws.Range("H63").Value = 4729.9
$ws.Range("J63").Value = 4954.5454
$ws.Range("L63").Value = 4954.5454
$ws.Range("N63").Value = -6326.5454
# Row 66
$ws.Range("H66").Value = 4729.9
$ws.Range("J66").Value = 4954.5454
$ws.Range("L66").Value = 24772.727
$ws.Range("N66").Value = -31636.727
# Row 116
$ws.Range("H116").Value = 1297.9697
$ws.Range("I116").Value = 1318.3182
$ws.Range("K116").Value = 1318.3182
$ws.Range("M116").Value = 975.6818000000001
# Row 132
$ws.Range("H132").Value = 14847.868
$ws.Range("I132").Value = 11385.241
$ws.Range("K132").Value = 34155.723
$ws.Range("M132").Value = -31625.723
# Row 136
$ws.Range("H136").Value = 4686.8213
$ws.Range("I136").Value = 4049.76
$ws.Range("K136").Value = 12149.28
$ws.Range("M136").Value = -9599.280000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1297.9697
$ws.Range("I3").Value = 1318.3182
$ws.Range("K3").Value = 1318.3182
$ws.Range("M3").Value = -1204.3182

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 315893.34
$ws.Range("I31").Value = 3001.5715
$ws.Range("J31").Value = 568613.6
$ws.Range("K31").Value = 3001.5715
$ws.Range("L31").Value = 568613.6
$ws.Range("M31").Value = -2706.5715
$ws.Range("N31").Value = -569203.6
# Row 34
$ws.Range("H34").Value = 315893.34
$ws.Range("I34").Value = 3001.5715
$ws.Range("J34").Value = 568613.6
$ws.Range("K34").Value = 3001.5715
$ws.Range("L34").Value = 568613.6
$ws.Range("M34").Value = -2799.5715
$ws.Range("N34").Value = -569017.6
# Row 41
$ws.Range("H41").Value = 37011.8
$ws.Range("J41").Value = 40000
$ws.Range("L41").Value = 40000
$ws.Range("N41").Value = -40856
# Row 58
$ws.Range("H58").Value = 2041
$ws.Range("I58").Value = 1709.3182
$ws.Range("K58").Value = 1709.3182
$ws.Range("M58").Value = -1506.3182
# Row 132
$ws.Range("H132").Value = 2440.111
$ws.Range("I132").Value = 2261.1853
$ws.Range("K132").Value = 6783.5559
$ws.Range("M132").Value = -4253.5559
# Row 134
$ws.Range("H134").Value = 1428.8
$ws.Range("I134").Value = 1382.7778
$ws.Range("J134").Value = 1497.8334
$ws.Range("K134").Value = 4148.3334
$ws.Range("L134").Value = 4493.5002
$ws.Range("M134").Value = -1613.3334
$ws.Range("N134").Value = -9563.5002
# Row 136
$ws.Range("H136").Value = 2041
$ws.Range("I136").Value = 1709.3182
$ws.Range("K136").Value = 5127.9546
$ws.Range("M136").Value = -2577.9546
# Row 141
$ws.Range("H141").Value = 246415.03
$ws.Range("J141").Value = 255680.7
$ws.Range("L141").Value = 255680.7
$ws.Range("N141").Value = -266040.7

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 10101368
$ws.Range("I2").Value = 75.72
$ws.Range("K2").Value = 454.32
$ws.Range("M2").Value = -341.32
# Row 7
$ws.Range("H7").Value = 3287.75
$ws.Range("I7").Value = 1175
$ws.Range("K7").Value = 3525
$ws.Range("M7").Value = -3413
# Row 12
$ws.Range("H12").Value = 107.2
$ws.Range("J12").Value = 149
$ws.Range("L12").Value = 447
$ws.Range("N12").Value = -793
# Row 23
$ws.Range("H23").Value = 580.8182
$ws.Range("I23").Value = 351.25
$ws.Range("J23").Value = 712
$ws.Range("K23").Value = 1053.75
$ws.Range("L23").Value = 2136
$ws.Range("M23").Value = -818.75
$ws.Range("N23").Value = -2606
# Row 34
$ws.Range("H34").Value = 1630
$ws.Range("I34").Value = 1630
$ws.Range("K34").Value = 4890
$ws.Range("M34").Value = -4806
# Row 39
$ws.Range("H39").Value = 7799.6
$ws.Range("J39").Value = 9499.5
$ws.Range("L39").Value = 28498.5
$ws.Range("N39").Value = -29086.5
# Row 86
$ws.Range("H86").Value = 2025.75
$ws.Range("J86").Value = 2025.75
$ws.Range("L86").Value = 6077.25
$ws.Range("N86").Value = -8449.25
# Row 89
$ws.Range("H89").Value = 2025.75
$ws.Range("J89").Value = 2025.75
$ws.Range("L89").Value = 18231.75
$ws.Range("N89").Value = -30087.75
# Row 113
$ws.Range("H113").Value = 47623636
$ws.Range("I113").Value = 4419.8
$ws.Range("K113").Value = 13259.4
$ws.Range("M113").Value = -11089.4
# Row 122
$ws.Range("H122").Value = 1366.3334
$ws.Range("J122").Value = 1374.5
$ws.Range("L122").Value = 12370.5
$ws.Range("N122").Value = -17270.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 17742.426
$ws.Range("I70").Value = 23262.8
$ws.Range("J70").Value = 8000.5884
$ws.Range("K70").Value = 23262.8
$ws.Range("L70").Value = 8000.5884
$ws.Range("M70").Value = -22992.8
$ws.Range("N70").Value = -8540.588400000001
# Row 73
$ws.Range("H73").Value = 17742.426
$ws.Range("I73").Value = 23262.8
$ws.Range("J73").Value = 8000.5884
$ws.Range("K73").Value = 23262.8
$ws.Range("L73").Value = 8000.5884
$ws.Range("M73").Value = -22326.8
$ws.Range("N73").Value = -9872.588400000001
# Row 102
$ws.Range("H102").Value = 2847.2812
$ws.Range("I102").Value = 2154.963
$ws.Range("J102").Value = 6585.8
$ws.Range("K102").Value = 2154.963
$ws.Range("L102").Value = 6585.8
$ws.Range("M102").Value = -532.9630000000002
$ws.Range("N102").Value = -9829.799999999999
# Row 122
$ws.Range("H122").Value = 4302.3335
$ws.Range("I122").Value = 4069.6667
$ws.Range("J122").Value = 5233
$ws.Range("K122").Value = 12209.0001
$ws.Range("L122").Value = 15699
$ws.Range("M122").Value = -9759.000100000001
$ws.Range("N122").Value = -20599

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 7626.8335
$ws.Range("I136").Value = 7699.2144
$ws.Range("K136").Value = 23097.6432
$ws.Range("M136").Value = -20547.6432
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 124995
$ws.Range("J94").Value = 124995
$ws.Range("L94").Value = 124995
$ws.Range("N94").Value = -126797
# Row 132
$ws.Range("H132").Value = 3289.2307
$ws.Range("I132").Value = 2240.8125
$ws.Range("K132").Value = 6722.4375
$ws.Range("M132").Value = -4192.4375
# Row 136
$ws.Range("H136").Value = 7495.5
$ws.Range("I136").Value = 6411.9287
$ws.Range("J136").Value = 9181.056
$ws.Range("K136").Value = 19235.7861
$ws.Range("L136").Value = 27543.168
$ws.Range("M136").Value = -16685.7861
$ws.Range("N136").Value = -32643.168
